$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two old rows (old sheet row 44 "A=44" at row 46, and old
# row "A=46" at row 48). Deleting bottom-most first keeps the earlier row
# index stable for the second delete. Every row below shifts up by one
# each time, pulling rows 47..52 up to 46..50.
$ws.Rows(48).Delete()
$ws.Rows(46).Delete()

# --- Append two brand-new rows (51 and 52) at the bottom of the table.
# Seed them from the last existing data row (50) first so they inherit its
# cell style (bold border s="1" on column A) instead of the unstyled
# default.
$ws.Range("A50").Copy($ws.Range("A51"))
$ws.Range("A50").Copy($ws.Range("A52"))

# Row 51: a "zero" placeholder row (summoner name blank, role SOLO).
$ws.Range("A51").Value = 51
$ws.Range("B51").Value = 0
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = ""
$ws.Range("H51").Value = "SOLO"
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0

# Row 52: new real data row.
$ws.Range("A52").Value = 52
$ws.Range("B52").Value = 2.719359230470012
$ws.Range("C52").Value = 3811.5
$ws.Range("D52").Value = 0.01770923299731227
$ws.Range("E52").Value = 28.5
$ws.Range("F52").Value = 279.5
$ws.Range("G52").Value = "Portgas D Åce "
$ws.Range("H52").Value = "SOLO"
$ws.Range("I52").Value = 0.1933071592924438
$ws.Range("J52").Value = 18
$ws.Range("K52").Value = 0.01137944447740923
